# FUNCTIONALITY: Delete - Tagged a test case and improved the position for
# automation. ListView - Finished automating the test suite for now.
# Read - Added the missing resource file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Data updates -----------------------------------------------------
# Create: test-case count adjustment (B2)
$ws.Range("B2").Value = 0

# ListView: improved automation coverage (B4) and finished status (D4)
$ws.Range("B4").Value = 7

# New annotation strings are entered in the same order in which they first
# appear in the finished workbook's shared-string table.
$ws.Range("E4").Value = "Contains two partially automated test cases."
$ws.Range("E5").Value = "Contains four partially automated test cases."

# Delete: tagged and repositioned for automation (D3/E3)
$ws.Range("D3").Value = "Ready to Write"
$ws.Range("E3").Value = "1 is Automateable, 2/3 are Automateable?"

# ListView: finished automating the test suite for now (D4)
$ws.Range("D4").Value = "Finished"

# --- Selection ----------------------------------------------------------
$ws.Range("E4").Select()
